$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for an "Unknown" vendor, with VendorID 0 and Name "Unknown"
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "Unknown"

# Move active selection to B5, matching the post-edit state
$ws.Range("B5").Select()
